$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Matière"
$ws.Range("B1").Value = "Salle"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Date de début"
$ws.Range("E1").Value = "Date de fin"
$ws.Range("F1").Value = "Durée"

# Row 2
$ws.Range("A2").Value = "M4203C-EVAL-2-EG-SB-BD"
$ws.Range("B2").Value = "A251 (85pl./160)"
$ws.Range("C2").Value = "\n\nG22-S4\nGRAZIANO EMMANUELLE\nBOUCHON STÉPHANIE\nDIARD BEN"
$ws.Range("D2").Value = 44645.5
$ws.Range("E2").Value = 44645.58333333334
$ws.Range("F2").Value = 0.08333333333333333

# Row 3
$ws.Range("A3").Value = "OS09A-EVAL-2"
$ws.Range("B3").Value = "D360 (16 pl./32)"
$ws.Range("C3").Value = "\n\nG22-S4\nDELL'AIERA MICHAEL\n(Exporté le:22/03/2022 18:44)"
$ws.Range("D3").Value = 44642.375
$ws.Range("E3").Value = 44642.45833333334
$ws.Range("F3").Value = 0.08333333333333333

# Row 4
$ws.Range("A4").Value = "OS03-EVAL-FJ"
$ws.Range("B4").Value = "D360 (16 pl./32)"
$ws.Range("C4").Value = "\n\nG22-S4\nJAZIRI FAOUZI\n(Exporté le:22/03/2022 18:44)\n"
$ws.Range("D4").Value = 44642.29166666666
$ws.Range("E4").Value = 44642.375
$ws.Range("F4").Value = 0.08333333333333333

# Row 5
$ws.Range("A5").Value = "Point stage"
$ws.Range("B5").Value = "A251 (85pl./160)"
$ws.Range("C5").Value = "\n\nINFO2-S4\nGRUSON NATHALIE\n(Exporté le:22/03/2022 18:44)\"
$ws.Range("D5").Value = 44644.47916666666
$ws.Range("E5").Value = 44644.52083333334
$ws.Range("F5").Value = 0.04166666666666666

# Row 6
$ws.Range("A6").Value = "TOEIC Examen DUETI"
$ws.Range("B6").Value = "E104 (66 pl./66)\,E102 (66 pl./66)\,E103 (50 pl./52)"
$ws.Range("C6").Value = "\n\nINFO2-S4\nBACHELET CAROLE\n(Exporté le:22/03/2022 18:44)\"
$ws.Range("D6").Value = 44644.52083333334
$ws.Range("E6").Value = 44644.64583333334
$ws.Range("F6").Value = 0.125

# Row 7
$ws.Range("A7").Value = "OS05B-EVAL-G22"
$ws.Range("B7").Value = "D250 (15 pl./27)\,D251 (11 pl./22)"
$ws.Range("C7").Value = "\n\nG22-S4\nDAMAS LUC\n(Exporté le:22/03/2022 18:44)\n"
$ws.Range("D7").Value = 44643.5
$ws.Range("E7").Value = 44643.66666666666
$ws.Range("F7").Value = 0.1666666666666667

# Number formats: columns D/E as dates, column F as time
$ws.Range("D2:E7").NumberFormat = "dd/mm/yy"
$ws.Range("F2:F7").NumberFormat = "hh:mm:ss"

# Page margins (inches -> points, 72 pts/in)
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
